# Minor fix in TSP.
# Update the "Fitness" column (C) values for rows 2-12 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3517, 3586, 3670, 3773, 3773, 4026, 4027, 4027, 4547, 4547, 4547)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
